# Update cryptos list prices (D) and 1h volume change (E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For D-column values that look like plain numbers (e.g. "1.002", "0.06650"),
# force the cell to Text format first so Excel does not auto-convert the
# string into a Number (which would silently drop significant trailing
# zeros / reformat the value). Cells whose new text is not numeric-looking
# (e.g. "27.089.60") do not need this and are left as plain text writes.

$ws.Range("D2").Value = "27.089.60"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "1.716.03"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.44"
$ws.Range("E5").Value = "  -6.00%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4716"
$ws.Range("E7").Value = "  +5.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3427"
$ws.Range("E8").Value = "  -3.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.11"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07268"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.043"
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.88"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.869"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").Value = "1.713.81"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.883"
$ws.Range("E16").Value = "  -5.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.19"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06355"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.51"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.622"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "27.130.65"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.121"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.17"
$ws.Range("E26").Value = "  -4.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.51"
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("D28").Value = "1.914.13"
$ws.Range("E28").Value = "  -3.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.091"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.53"
$ws.Range("E30").Value = "  -4.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.015"
$ws.Range("E31").Value = "  -8.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09151"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.594"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.323"
$ws.Range("E34").Value = "  -5.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02203"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05826"
$ws.Range("E36").Value = "  -4.39%  "
$ws.Range("E37").Value = "  -7.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.1997"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.736"
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.394"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5891"
$ws.Range("E41").Value = "  -7.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.117"
$ws.Range("E42").Value = "  -6.07%  "
$ws.Range("E43").Value = "  -5.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.55"
$ws.Range("E44").Value = "  -5.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5651"
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.563"
$ws.Range("E46").Value = "  -4.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.33"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.842"
$ws.Range("E48").Value = "  -5.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06650"
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("E51").Value = "  +0.08%  "

Write-Host "Updated cryptos prices/volumes"
